$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before the current row 12 (shifts old rows 12-19 down to 15-22)
$ws.Range("A12:M14").EntireRow.Insert()

# Fill the newly inserted row 12 with the new "NOP" instruction data
$ws.Range("A12").Value = "NOP"
$ws.Range("B12").Value = "``1111"
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0

# Update selected cell to match the saved view state
$ws.Range("A13").Select()
